# B6-PowerPoint.pptx edit
#
# 1) Re-apply the table style used by the three tables (slides 14-16) so
#    they pick up the built-in "{AFA9816F-A0C0-4A95-9DC4-091501C23DF4}"
#    gallery style instead of the custom "Table_0" style that was applied
#    before.
# 2) Re-colour the deck's theme (ppt/theme/theme1.xml, the theme used by
#    the slide master) so it matches the "Office" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Tables: swap the applied table style on every table in the deck.
# ---------------------------------------------------------------------
$newTableStyleId = "{AFA9816F-A0C0-4A95-9DC4-091501C23DF4}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Theme colours: apply the "Office" colour scheme (dk1/lt1/dk2/lt2/
#    accent1-6/hlink/folHlink) to the presentation's theme.
# ---------------------------------------------------------------------
function HexToOleRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $themeColors.Count; $k++) {
    $themeColors.Item($k).RGB = HexToOleRgb $officeColors[$k - 1]
}
